# MSME Country Indicators - Luxembourg Summary
# Source section (rows ~70-80 on "Summary") is restructured:
#  - A new blank line is inserted above the "SBS Main Indicators..." source line.
#  - The Eurostat data-source hyperlink is removed (cell keeps its text, loses the
#    live hyperlink + blue/underline styling) and moves down one row, with an
#    extra blank line inserted above it (reusing the old blank line that used to
#    sit right after it).
#  - The old Luxembourg statistics-portal citation text is replaced with a
#    citation for the Luxembourg "Règlement grand-ducal" (MSME definition) text.
#  - The old "Sructural Business Statistics - Eurostat..." citation line is
#    replaced by a second "SBS Eurostat" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hyperlink that lives on A72 ("http://epp.eurostat...") before the
# row-insert shifts the sheet contents down, so we don't have to chase it.
$ws.Range("A72").Hyperlinks.Delete()

# Insert one new row above row 71 ("SBS Main Indicators..."); this pushes
# everything from the old row 71 down through row 79 down by one row, so the
# old row 73 (blank) / row 72 (url) pairing now sits at rows 74 / 73.
$ws.Rows("71").Insert()

# The blank row that used to directly follow the url line (old A73, now A74)
# should instead sit directly above the url line (now A74), i.e. swap the
# "blank" and "url" cell contents that the insert left in rows 73/74.
$ws.Range("A73").Value = ""
$ws.Range("A73").Font.Italic = $true

$ws.Range("A74").Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"

# Replace the Statistics Portal citation text with the Règlement grand-ducal one.
$ws.Range("A78").Value = "Règlement grand-ducal, Mémorial A n° 38 de 2005, Définition des micro, petites et moyennes entreprises available at http://www.legilux.public.lu/leg/a/archives/2005/0038/index.html"

# Replace the old "Sructural Business Statistics..." citation with "SBS Eurostat".
$ws.Range("A80").Value = "SBS Eurostat"
